$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductionPlan")

# Update existing row 2 (MAT_B -> MAT_A) with new quantities
$ws.Range("A2").Value = "MAT_A"
$ws.Range("C2").Value = "LINE_A"
$ws.Range("G2").Value = 870
$ws.Range("H2").Value = 870
$ws.Range("J2").Value = 827

# Add new row 3 with the original MAT_B / LINE_B data and new quantities
$ws.Range("A3").Value = "MAT_B"
$ws.Range("B3").Value = "PLANT_001"
$ws.Range("C3").Value = "LINE_B"
$ws.Range("D3").Value = $ws.Range("D2").Value2
$ws.Range("E3").Value = $ws.Range("E2").Value2
$ws.Range("F3").Value = $ws.Range("F2").Value2
$ws.Range("D3:F3").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("G3").Value = 120
$ws.Range("H3").Value = 120
$ws.Range("I2").Copy($ws.Range("I3"))
$ws.Range("J3").Value = 106
